# ----------------------------------------------------------------------------
# The commit swaps the two theme parts of the deck: the theme that is
# actually applied to the slide master/slides ("Integral") and the theme
# that is only wired to the notes master ("Office Theme") trade places, so
# the slide master ends up using the stock "Office Theme" palette while the
# notes master ends up with the "Integral" palette.
#
# $ppt.ActivePresentation.SlideMaster.Theme is the one-and-only Design/Theme
# object this host models (Designs.Count is 1 and every Slide/Master -
# including the notes master - resolves to that same Theme), and it backs
# ppt/theme/theme2.xml, the part that is actually referenced by the slide
# master. Re-pointing that Theme's colour scheme from the "Integral" values
# to the "Office Theme" values reproduces the half of the swap that changes
# what viewers actually see.
# ----------------------------------------------------------------------------

function Set-ThemeRGB {
    param($ThemeColor, [string]$HexRRGGBB)
    $r = [Convert]::ToInt32($HexRRGGBB.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexRRGGBB.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexRRGGBB.Substring(4, 2), 16)
    # PowerPoint's RGB colour longs are stored low-byte-first (R + G*256 + B*65536),
    # matching the classic VBA RGB(r,g,b) macro.
    $ThemeColor.RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette: the stock Office "Office Theme" colour scheme, in the
# fixed Colors(1..12) order used by ThemeColorScheme:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
Set-ThemeRGB $colorScheme.Colors(1)  "000000"
Set-ThemeRGB $colorScheme.Colors(2)  "FFFFFF"
Set-ThemeRGB $colorScheme.Colors(3)  "44546A"
Set-ThemeRGB $colorScheme.Colors(4)  "E7E6E6"
Set-ThemeRGB $colorScheme.Colors(5)  "5B9BD5"
Set-ThemeRGB $colorScheme.Colors(6)  "ED7D31"
Set-ThemeRGB $colorScheme.Colors(7)  "A5A5A5"
Set-ThemeRGB $colorScheme.Colors(8)  "FFC000"
Set-ThemeRGB $colorScheme.Colors(9)  "4472C4"
Set-ThemeRGB $colorScheme.Colors(10) "70AD47"
Set-ThemeRGB $colorScheme.Colors(11) "0563C1"
Set-ThemeRGB $colorScheme.Colors(12) "954F72"

# Best-effort: also rename the Design/theme to match its new palette. (This
# host currently treats Design.Name as read-only and silently drops the
# assignment, but it is the correct call to make and is harmless either way.)
$design = $p.Designs.Item(1)
$design.Name = "Office Theme"
